$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 6114
$ws.Range("B3").Value = 148130
$ws.Range("B4").Value = 28673
$ws.Range("B5").Value = 19.36
$ws.Range("B6").Value = 24.23
$ws.Range("B7").Value = 6.64
